# edit.ps1 - reproduces:
#   1) the table style-id change on the B1/B2 slide's table
#      ({1827EAF6-...} -> {CBD470F4-...})
#   2) the theme colour swap between the two theme parts: the
#      presentation's (only reachable/writable) theme - "Integral" /
#      "Red Violet" - is recoloured to the stock "Office" palette that
#      used to live in the sibling (orphaned, notes-master-only) theme
#      part.

$p = $ppt.ActivePresentation

# --- 1) Table style id -------------------------------------------------
$s   = $p.Slides.Item(5)
$shp = $s.Shapes.Item(2)
if ($shp.HasTable) {
    $tbl = $shp.Table
    $tbl.ApplyStyle("{CBD470F4-34C0-49FC-8340-6F09E9F889FC}")
}

# --- 2) Theme colours ----------------------------------------------------
# Target palette ("Office") expressed as VBA-style BGR-packed longs
# (OLE_COLOR / RGB() layout: 0x00BBGGRR) so that Colors(i).RGB round-trips
# to the correct <a:srgbClr val="RRGGBB"/> on save.
$tcs = $s.ThemeColorScheme

$tcs.Colors(1).RGB  = 0x000000   # dk1      -> 000000
$tcs.Colors(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$tcs.Colors(3).RGB  = 0x6A5444   # dk2      -> 44546A
$tcs.Colors(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$tcs.Colors(6).RGB  = 0x317DED   # accent2  -> ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$tcs.Colors(9).RGB  = 0xC47244   # accent5  -> 4472C4
$tcs.Colors(10).RGB = 0x47AD70   # accent6  -> 70AD47
$tcs.Colors(11).RGB = 0xC16305   # hlink    -> 0563C1
$tcs.Colors(12).RGB = 0x724F95   # folHlink -> 954F72

# Best-effort rename (a no-op on this host since ThemeColorScheme.Name is
# read-only here, but harmless and matches the authored clrScheme name if
# the host ever supports it).
try { $tcs.Name = "Office" } catch {}
